# Applies the cryptos.xlsx data refresh described in the commit diff.
# Only cell VALUES change (no styles/formats); B/C/D/E text cells are
# updated to match the new scrape snapshot, including a couple of row
# swaps (Litecoin/Polygon, TrustWalletToken/Quant).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.343.94'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.879.83'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Formula = '="0.7107"'
$ws.Range('D5').Copy() | Out-Null
$ws.Range('D5').PasteSpecial(-4163) | Out-Null
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').Formula = '="242.54"'
$ws.Range('D6').Copy() | Out-Null
$ws.Range('D6').PasteSpecial(-4163) | Out-Null
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Formula = '="0.08020"'
$ws.Range('D8').Copy() | Out-Null
$ws.Range('D8').PasteSpecial(-4163) | Out-Null
$ws.Range('E8').Value = '  +3.05%  '
$ws.Range('D9').Formula = '="0.3138"'
$ws.Range('D9').Copy() | Out-Null
$ws.Range('D9').PasteSpecial(-4163) | Out-Null
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').Formula = '="0.08329"'
$ws.Range('D11').Copy() | Out-Null
$ws.Range('D11').PasteSpecial(-4163) | Out-Null
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('D12').Value = '1.891.74'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').Formula = '="5.268"'
$ws.Range('D13').Copy() | Out-Null
$ws.Range('D13').PasteSpecial(-4163) | Out-Null
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Formula = '="0.7180"'
$ws.Range('D14').Copy() | Out-Null
$ws.Range('D14').PasteSpecial(-4163) | Out-Null
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Formula = '="94.35"'
$ws.Range('D15').Copy() | Out-Null
$ws.Range('D15').PasteSpecial(-4163) | Out-Null
$ws.Range('E15').Value = '  +3.39%  '
$ws.Range('D16').Formula = '="6.358"'
$ws.Range('D16').Copy() | Out-Null
$ws.Range('D16').PasteSpecial(-4163) | Out-Null
$ws.Range('E16').Value = '  +5.09%  '
$ws.Range('D17').Formula = '="0.000008709"'
$ws.Range('D17').Copy() | Out-Null
$ws.Range('D17').PasteSpecial(-4163) | Out-Null
$ws.Range('E17').Value = '  +5.31%  '
$ws.Range('D18').Value = '29.358.81'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('D19').Formula = '="242.87"'
$ws.Range('D19').Copy() | Out-Null
$ws.Range('D19').PasteSpecial(-4163) | Out-Null
$ws.Range('E19').Value = '  +0.58%  '
$ws.Range('D20').Value = '2.150.99'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').Formula = '="13.33"'
$ws.Range('D21').Copy() | Out-Null
$ws.Range('D21').PasteSpecial(-4163) | Out-Null
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').Formula = '="7.846"'
$ws.Range('D23').Copy() | Out-Null
$ws.Range('D23').PasteSpecial(-4163) | Out-Null
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('D24').Formula = '="1.002"'
$ws.Range('D24').Copy() | Out-Null
$ws.Range('D24').PasteSpecial(-4163) | Out-Null
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').Formula = '="0.1570"'
$ws.Range('D25').Copy() | Out-Null
$ws.Range('D25').PasteSpecial(-4163) | Out-Null
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').Formula = '="163.70"'
$ws.Range('D26').Copy() | Out-Null
$ws.Range('D26').PasteSpecial(-4163) | Out-Null
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Formula = '="9.078"'
$ws.Range('D27').Copy() | Out-Null
$ws.Range('D27').PasteSpecial(-4163) | Out-Null
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Formula = '="18.60"'
$ws.Range('D28').Copy() | Out-Null
$ws.Range('D28').PasteSpecial(-4163) | Out-Null
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').Formula = '="4.435"'
$ws.Range('D30').Copy() | Out-Null
$ws.Range('D30').PasteSpecial(-4163) | Out-Null
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').Formula = '="4.356"'
$ws.Range('D31').Copy() | Out-Null
$ws.Range('D31').PasteSpecial(-4163) | Out-Null
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').Formula = '="1.204"'
$ws.Range('D32').Copy() | Out-Null
$ws.Range('D32').PasteSpecial(-4163) | Out-Null
$ws.Range('E32').Value = '  -6.49%  '
$ws.Range('D33').Formula = '="0.05392"'
$ws.Range('D33').Copy() | Out-Null
$ws.Range('D33').PasteSpecial(-4163) | Out-Null
$ws.Range('E33').Value = '  +1.57%  '
$ws.Range('D34').Formula = '="1.944"'
$ws.Range('D34').Copy() | Out-Null
$ws.Range('D34').PasteSpecial(-4163) | Out-Null
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Formula = '="0.7793"'
$ws.Range('D35').Copy() | Out-Null
$ws.Range('D35').PasteSpecial(-4163) | Out-Null
$ws.Range('E35').Value = '  +4.14%  '
$ws.Range('D36').Formula = '="1.178"'
$ws.Range('D36').Copy() | Out-Null
$ws.Range('D36').PasteSpecial(-4163) | Out-Null
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').Formula = '="2.689"'
$ws.Range('D37').Copy() | Out-Null
$ws.Range('D37').PasteSpecial(-4163) | Out-Null
$ws.Range('D38').Formula = '="0.01888"'
$ws.Range('D38').Copy() | Out-Null
$ws.Range('D38').PasteSpecial(-4163) | Out-Null
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('D39').Value = '1.271.36'
$ws.Range('E39').Value = '  +5.16%  '
$ws.Range('D40').Formula = '="2.748"'
$ws.Range('D40').Copy() | Out-Null
$ws.Range('D40').PasteSpecial(-4163) | Out-Null
$ws.Range('E40').Value = '  +0.93%  '
$ws.Range('E41').Value = '  +1.58%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Formula = '="113.90"'
$ws.Range('D42').Copy() | Out-Null
$ws.Range('D42').PasteSpecial(-4163) | Out-Null
$ws.Range('E42').Value = '  +3.91%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Formula = '="0.9190"'
$ws.Range('D43').Copy() | Out-Null
$ws.Range('D43').PasteSpecial(-4163) | Out-Null
$ws.Range('E43').Value = '  +3.46%  '
$ws.Range('D44').Formula = '="74.57"'
$ws.Range('D44').Copy() | Out-Null
$ws.Range('D44').PasteSpecial(-4163) | Out-Null
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').Value = '2.036.21'
$ws.Range('E46').Value = '  +0.84%  '
$ws.Range('E47').Value = '  +3.76%  '
$ws.Range('D48').Formula = '="1.807"'
$ws.Range('D48').Copy() | Out-Null
$ws.Range('D48').PasteSpecial(-4163) | Out-Null
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('D49').Formula = '="0.5221"'
$ws.Range('D49').Copy() | Out-Null
$ws.Range('D49').PasteSpecial(-4163) | Out-Null
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').Formula = '="9.554"'
$ws.Range('D50').Copy() | Out-Null
$ws.Range('D50').PasteSpecial(-4163) | Out-Null
$ws.Range('E50').Value = '  +1.77%  '
$ws.Range('D51').Formula = '="0.4384"'
$ws.Range('D51').Copy() | Out-Null
$ws.Range('D51').PasteSpecial(-4163) | Out-Null
$ws.Range('E51').Value = '  +1.39%  '

$excel.CutCopyMode = 0
